# "Updates for fall 2018" - refresh the Excel/QlikView/R comparison data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3's task was renamed: "update data" -> "separate data & analysis"
$ws.Range("B3").Value = "separate data & analysis"

# Excel's rating for "reproducible analysis" (row 4, col R/E) improved: good -> best
$ws.Range("E4").Value = "best"

# Row 5's task was renamed: "workflow flexibility" -> "data exploration"
$ws.Range("B5").Value = "data exploration"

# Excel's rating for "complex plotting" (row 11, col R/E) improved: good -> best
$ws.Range("E11").Value = "best"

# Drop the stray leftover note cell below the table ("adsf")
$ws.Range("C16").ClearContents()

# Move the selection to reflect where the edits were made
$ws.Range("D5").Select()
